$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (file, name) values for rows 5..19, matching the target diff.
$data = @(
    @("070cef68-be1f-4e22-a9c8-d48853145b94_fila_14.png", "Angel Isaac Alvarez Quiñonez"),
    @("20a84539-ad79-4116-9810-100a1db56478_fila_15.png", "Orlando Mauricio Guevara"),
    @("250cbdd7-aea0-4d12-9022-ed21042ac5dc_fila_8.png", "David Orlando Mena Valverd"),
    @("2b23fdb2-bdc6-4ec4-92b6-3579de34c8e1_fila_4.png", "Yadder Fernando Torres"),
    @("3cd051f9-b9f8-4e06-8442-37f522ee2e44_fila_1.png", "Hotep Antonio Ruiz Lezama"),
    @("5850187f-9762-47cf-8cdd-f1c7f9202558_fila_10.png", "Abraham Silva Ampre"),
    @("691a6f4e-91d7-4318-ab00-453157c744c5_fila_6.png", "Ronier Jose Rivera"),
    @("6c09a82a-1687-4c4b-a651-0ae592b46dcd_fila_9.png", "Roman Alfonso Grios Boza"),
    @("7d239cd3-2e99-42be-a368-4a7ad19da24d_fila_11.png", "Eduardo Domingo Zeledon Merca"),
    @("9c9f8199-36dc-43cd-b462-59ebef39d290_fila_13.png", "Marlon Josue Gonzales Cano"),
    @("a21ce003-bcbd-4e0e-83bd-a34d92403ed7_fila_2.png", "Isabella Dompe Estrada"),
    @("a635feb3-33c4-4ee5-842a-7c73cccfc65e_fila_5.png", "Erick Espinoza"),
    @("bb2ca6a3-5dd2-415f-85af-6528d1e39552_fila_3.png", "Bryan Alexander Cano"),
    @("cf8a472c-3b8e-4162-8e62-3ec0022c9130_fila_12.png", "José Danilo Suárez"),
    @("e9442c21-0052-45c3-ab87-27b2cbc45ef2_fila_7.png", "Cristina Jozabed Carvajal")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
